$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (border + wrap text + row height) from row 10 down to the new rows 11-13
$ws.Range("A10:AY10").Copy()
$ws.Range("A11:AY13").PasteSpecial(-4122)
$ws.Range("11:13").RowHeight = 32

$arr11 = New-Object 'object[,]' 1,49
$arr11[0,0] = "2015C"
$arr11[0,1] = "ENGR"
$arr11[0,2] = 100
$arr11[0,3] = 500
$arr11[0,4] = "Scott Deaton"
$arr11[0,5] = 40
$arr11[0,6] = "LEC"
$arr11[0,7] = 7
$arr11[0,8] = 36.840000000000003
$arr11[0,9] = 3.3
$arr11[0,10] = 0.82
$arr11[0,11] = 3.14
$arr11[0,12] = 1.069
$arr11[0,13] = 3.43
$arr11[0,14] = 1.397
$arr11[0,15] = 3
$arr11[0,16] = 0.81599999999999995
$arr11[0,17] = 3.71
$arr11[0,18] = 0.75600000000000001
$arr11[0,19] = 3.71
$arr11[0,20] = 1.113
$arr11[0,21] = 3.71
$arr11[0,22] = 0.48799999999999999
$arr11[0,23] = 3.71
$arr11[0,24] = 0.48799999999999999
$arr11[0,25] = 3.71
$arr11[0,26] = 0.95099999999999996
$arr11[0,27] = 3.14
$arr11[0,28] = 0.69
$arr11[0,29] = 2.86
$arr11[0,30] = 0.378
$arr11[0,31] = 3.29
$arr11[0,32] = 0.75600000000000001
$arr11[0,33] = 3
$arr11[0,34] = 0.57699999999999996
$arr11[0,35] = 3.14
$arr11[0,36] = 0.9
$arr11[0,37] = 2.86
$arr11[0,38] = 1.069
$arr11[0,39] = 2.86
$arr11[0,40] = 0.378
$arr11[0,41] = 3.57
$arr11[0,42] = 0.78700000000000003
$arr11[0,43] = 3.29
$arr11[0,44] = 0.75600000000000001
$arr11[0,45] = 3.29
$arr11[0,46] = 0.48799999999999999
$arr11[0,47] = 3.33
$arr11[0,48] = 0.81599999999999995
$ws.Range("A11:AW11").Value = $arr11

$arr12 = New-Object 'object[,]' 1,49
$arr12[0,0] = "2015C"
$arr12[0,1] = "ENGR"
$arr12[0,2] = 110
$arr12[0,3] = 500
$arr12[0,4] = "Not CSCE"
$arr12[0,5] = 40
$arr12[0,6] = "LEC"
$arr12[0,7] = 7
$arr12[0,8] = 36.840000000000003
$arr12[0,9] = 3.3
$arr12[0,10] = 0.82
$arr12[0,11] = 3.14
$arr12[0,12] = 1.069
$arr12[0,13] = 3.43
$arr12[0,14] = 1.397
$arr12[0,15] = 3
$arr12[0,16] = 0.81599999999999995
$arr12[0,17] = 3.71
$arr12[0,18] = 0.75600000000000001
$arr12[0,19] = 3.71
$arr12[0,20] = 1.113
$arr12[0,21] = 3.71
$arr12[0,22] = 0.48799999999999999
$arr12[0,23] = 3.71
$arr12[0,24] = 0.48799999999999999
$arr12[0,25] = 3.71
$arr12[0,26] = 0.95099999999999996
$arr12[0,27] = 3.14
$arr12[0,28] = 0.69
$arr12[0,29] = 2.86
$arr12[0,30] = 0.378
$arr12[0,31] = 3.29
$arr12[0,32] = 0.75600000000000001
$arr12[0,33] = 3
$arr12[0,34] = 0.57699999999999996
$arr12[0,35] = 3.14
$arr12[0,36] = 0.9
$arr12[0,37] = 2.86
$arr12[0,38] = 1.069
$arr12[0,39] = 2.86
$arr12[0,40] = 0.378
$arr12[0,41] = 3.57
$arr12[0,42] = 0.78700000000000003
$arr12[0,43] = 3.29
$arr12[0,44] = 0.75600000000000001
$arr12[0,45] = 3.29
$arr12[0,46] = 0.48799999999999999
$arr12[0,47] = 3.33
$arr12[0,48] = 0.81599999999999995
$ws.Range("A12:AW12").Value = $arr12

$arr13 = New-Object 'object[,]' 1,49
$arr13[0,0] = "2015C"
$arr13[0,1] = "PETE"
$arr13[0,2] = 100
$arr13[0,3] = 500
$arr13[0,4] = "Not CSCE"
$arr13[0,5] = 40
$arr13[0,6] = "LEC"
$arr13[0,7] = 7
$arr13[0,8] = 36.840000000000003
$arr13[0,9] = 3.3
$arr13[0,10] = 0.82
$arr13[0,11] = 3.14
$arr13[0,12] = 1.069
$arr13[0,13] = 3.43
$arr13[0,14] = 1.397
$arr13[0,15] = 3
$arr13[0,16] = 0.81599999999999995
$arr13[0,17] = 3.71
$arr13[0,18] = 0.75600000000000001
$arr13[0,19] = 3.71
$arr13[0,20] = 1.113
$arr13[0,21] = 3.71
$arr13[0,22] = 0.48799999999999999
$arr13[0,23] = 3.71
$arr13[0,24] = 0.48799999999999999
$arr13[0,25] = 3.71
$arr13[0,26] = 0.95099999999999996
$arr13[0,27] = 3.14
$arr13[0,28] = 0.69
$arr13[0,29] = 2.86
$arr13[0,30] = 0.378
$arr13[0,31] = 3.29
$arr13[0,32] = 0.75600000000000001
$arr13[0,33] = 3
$arr13[0,34] = 0.57699999999999996
$arr13[0,35] = 3.14
$arr13[0,36] = 0.9
$arr13[0,37] = 2.86
$arr13[0,38] = 1.069
$arr13[0,39] = 2.86
$arr13[0,40] = 0.378
$arr13[0,41] = 3.57
$arr13[0,42] = 0.78700000000000003
$arr13[0,43] = 3.29
$arr13[0,44] = 0.75600000000000001
$arr13[0,45] = 3.29
$arr13[0,46] = 0.48799999999999999
$arr13[0,47] = 3.33
$arr13[0,48] = 0.81599999999999995
$ws.Range("A13:AW13").Value = $arr13

# AX/AY columns stay blank (already created by the format paste above)

# Update active selection to match the end-state of editing (B13)
$ws.Range("B13").Select()
